$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: A22 gets a date value (2020-03-23 -> serial 43913) ---
# Copy the date-format style from A20 (style index 4: numFmtId 14, centered)
# so we reuse the existing style instead of minting a new one.
$ws.Range("A20").Copy()
$ws.Range("A22").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A22").Value = 43913

# --- Row 23: new journal entry ---
$ws.Range("A23").Value = "30.02.2020"
$ws.Range("B23").Value = 8
$ws.Range("D23").Value = "4h"
$ws.Range("C23").Value = "Création de grille aléatoire du jeu"

# F23 changes style (s=1 -> s=5, wrap/center like E23) as well as getting text.
$ws.Range("E23").Copy()
$ws.Range("F23").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F23").Value = "Les bateaux des dernières cases appart le bateau qui a juste une case, qui apparaissait a la fin du tableau s'écrivait mais n'était pas visible. J'ai trouvé une solution qui consite a limiter jusqu'à ou les bateaux peuvent aller comme ça ne dépasse pas de la grille (Je m'excuse si mon explication n'est pas très compréhensible)"

# Row 23 grows tall to fit the wrapped comment.
$ws.Rows.Item(23).RowHeight = 158.4

# --- View state: selection moves to F40 after scrolling to A22 ---
$ws.Range("A22").Select() | Out-Null
$ws.Range("F40").Select() | Out-Null
